$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text (conversion rates) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 2.33 = 8771.99 pesos"), "✅ 1000 Bs = 2.34 = 8822.45 pesos"
$text = $text -replace [regex]::Escape("✅ 8771.99 pesos = 2.31 = 960.3 Bs"), "✅ 8822.45 pesos = 2.32 = 957.8 Bs"
$wsHoja1.Range("A1").Value2 = $text

# --- Update tasas sheet numeric cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 428
$wsTasas.Range("O10").Value = 3776.01
$wsTasas.Range("N12").Value = 3795
$wsTasas.Range("O12").Value = 412
